$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregado de etiqueta 303-HE-PECHO 3A: se carga el numero de etiqueta
# y se marca como traducida (OK) en la fila 11.
$ws.Range("M11").Value = 303
$ws.Range("N11").Value = "OK"

# La fila 11 crece para acomodar el contenido nuevo.
$ws.Rows(11).RowHeight = 84.75

# La celda N12 (ya en estado OK) se resalta con subrayado para destacar
# la nueva etiqueta agregada.
$ws.Range("N12").Font.Underline = 2

# Selección final tal cual quedó en el archivo original.
$null = $ws.Range("N12").Select()
